$p = $ppt.ActivePresentation

# The sentence lives on the "summary" slide, in the body placeholder that
# lists the retro bullet points ("What went well / What went bad / What
# would I change ..."). Locate it defensively by scanning all slides/shapes
# instead of hard-coding indices.
$enDash = [char]0x2013
$target = "What would I change"

$slide = $null
$shape = $null
foreach ($s in $p.Slides) {
    foreach ($sh in $s.Shapes) {
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text.Contains($target)) {
                $slide = $s
                $shape = $sh
            }
        }
    }
}

$tr = $shape.TextFrame.TextRange

# Find the paragraph that contains the "What would I change" bullet. (The
# reported .Paragraphs().Count is unreliable in this host, so walk a
# generous range instead -- out-of-range Paragraphs(i,1) calls simply come
# back empty rather than throwing.)
$para = $null
for ($i = 1; $i -le 50; $i++) {
    $candidate = $tr.Paragraphs($i, 1)
    if ($candidate.Text.Length -gt 0 -and $candidate.Text.StartsWith($target)) {
        $para = $candidate
    }
}

$fullText = $para.Text
$oldTail = "additional methods for the CRUD such as delete-all which was planned"
$splitIndex = $fullText.IndexOf($oldTail)

# First run keeps everything up to (and including) "possibly ", plus the
# new word "some " that now belongs with it.
$headRange = $para.Characters(1, $splitIndex)

# Second run becomes the brand-new sentence fragment.
$tailRange = $para.Characters($splitIndex + 1, $fullText.Length - $splitIndex)

# Replace the tail first (while headRange's span is still untouched),
# then extend the head text so it ends in "...possibly some ".
$tailRange.Text = "additional modals too"
$headRange.Text = "What would I change " + $enDash + " add some animation to the front end of the application, possibly some "
